$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")
$ws.Range("A1").Value = "test"
